$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot "before" values for the columns that move (D,H,I,J,K,L,M,O,P)
# across the affected rows, then reassign per the row permutation.
# (Use .Value2 -- .Value getter does not resolve reliably in this host.)
$snap = @{}
$snap[2] = @{
    D = $ws.Range("D2").Value2
    H = $ws.Range("H2").Value2
    I = $ws.Range("I2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    O = $ws.Range("O2").Value2
    P = $ws.Range("P2").Value2
}
$snap[3] = @{
    D = $ws.Range("D3").Value2
    H = $ws.Range("H3").Value2
    I = $ws.Range("I3").Value2
    J = $ws.Range("J3").Value2
    K = $ws.Range("K3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    O = $ws.Range("O3").Value2
    P = $ws.Range("P3").Value2
}
$snap[4] = @{
    D = $ws.Range("D4").Value2
    H = $ws.Range("H4").Value2
    I = $ws.Range("I4").Value2
    J = $ws.Range("J4").Value2
    K = $ws.Range("K4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    O = $ws.Range("O4").Value2
    P = $ws.Range("P4").Value2
}
$snap[5] = @{
    D = $ws.Range("D5").Value2
    H = $ws.Range("H5").Value2
    I = $ws.Range("I5").Value2
    J = $ws.Range("J5").Value2
    K = $ws.Range("K5").Value2
    L = $ws.Range("L5").Value2
    M = $ws.Range("M5").Value2
    O = $ws.Range("O5").Value2
    P = $ws.Range("P5").Value2
}
$snap[6] = @{
    D = $ws.Range("D6").Value2
    H = $ws.Range("H6").Value2
    I = $ws.Range("I6").Value2
    J = $ws.Range("J6").Value2
    K = $ws.Range("K6").Value2
    L = $ws.Range("L6").Value2
    M = $ws.Range("M6").Value2
    O = $ws.Range("O6").Value2
    P = $ws.Range("P6").Value2
}
$snap[8] = @{
    D = $ws.Range("D8").Value2
    H = $ws.Range("H8").Value2
    I = $ws.Range("I8").Value2
    J = $ws.Range("J8").Value2
    K = $ws.Range("K8").Value2
    L = $ws.Range("L8").Value2
    M = $ws.Range("M8").Value2
    O = $ws.Range("O8").Value2
    P = $ws.Range("P8").Value2
}
$snap[10] = @{
    D = $ws.Range("D10").Value2
    H = $ws.Range("H10").Value2
    I = $ws.Range("I10").Value2
    J = $ws.Range("J10").Value2
    K = $ws.Range("K10").Value2
    L = $ws.Range("L10").Value2
    M = $ws.Range("M10").Value2
    O = $ws.Range("O10").Value2
    P = $ws.Range("P10").Value2
}
$snap[11] = @{
    D = $ws.Range("D11").Value2
    H = $ws.Range("H11").Value2
    I = $ws.Range("I11").Value2
    J = $ws.Range("J11").Value2
    K = $ws.Range("K11").Value2
    L = $ws.Range("L11").Value2
    M = $ws.Range("M11").Value2
    O = $ws.Range("O11").Value2
    P = $ws.Range("P11").Value2
}
$snap[12] = @{
    D = $ws.Range("D12").Value2
    H = $ws.Range("H12").Value2
    I = $ws.Range("I12").Value2
    J = $ws.Range("J12").Value2
    K = $ws.Range("K12").Value2
    L = $ws.Range("L12").Value2
    M = $ws.Range("M12").Value2
    O = $ws.Range("O12").Value2
    P = $ws.Range("P12").Value2
}
$snap[13] = @{
    D = $ws.Range("D13").Value2
    H = $ws.Range("H13").Value2
    I = $ws.Range("I13").Value2
    J = $ws.Range("J13").Value2
    K = $ws.Range("K13").Value2
    L = $ws.Range("L13").Value2
    M = $ws.Range("M13").Value2
    O = $ws.Range("O13").Value2
    P = $ws.Range("P13").Value2
}
$snap[15] = @{
    D = $ws.Range("D15").Value2
    H = $ws.Range("H15").Value2
    I = $ws.Range("I15").Value2
    J = $ws.Range("J15").Value2
    K = $ws.Range("K15").Value2
    L = $ws.Range("L15").Value2
    M = $ws.Range("M15").Value2
    O = $ws.Range("O15").Value2
    P = $ws.Range("P15").Value2
}
$snap[16] = @{
    D = $ws.Range("D16").Value2
    H = $ws.Range("H16").Value2
    I = $ws.Range("I16").Value2
    J = $ws.Range("J16").Value2
    K = $ws.Range("K16").Value2
    L = $ws.Range("L16").Value2
    M = $ws.Range("M16").Value2
    O = $ws.Range("O16").Value2
    P = $ws.Range("P16").Value2
}
$snap[17] = @{
    D = $ws.Range("D17").Value2
    H = $ws.Range("H17").Value2
    I = $ws.Range("I17").Value2
    J = $ws.Range("J17").Value2
    K = $ws.Range("K17").Value2
    L = $ws.Range("L17").Value2
    M = $ws.Range("M17").Value2
    O = $ws.Range("O17").Value2
    P = $ws.Range("P17").Value2
}

# Apply permuted values: row r gets the pre-edit values of row mapping[r]
$ws.Range("D2").Value2 = $snap[5].D
$ws.Range("H2").Value2 = $snap[5].H
$ws.Range("I2").Value2 = $snap[5].I
$ws.Range("J2").Value2 = $snap[5].J
$ws.Range("K2").Value2 = $snap[5].K
$ws.Range("L2").Value2 = $snap[5].L
$ws.Range("M2").Value2 = $snap[5].M
$ws.Range("O2").Value2 = $snap[5].O
$ws.Range("P2").Value2 = $snap[5].P

$ws.Range("D3").Value2 = $snap[8].D
$ws.Range("H3").Value2 = $snap[8].H
$ws.Range("I3").Value2 = $snap[8].I
$ws.Range("J3").Value2 = $snap[8].J
$ws.Range("K3").Value2 = $snap[8].K
$ws.Range("L3").Value2 = $snap[8].L
$ws.Range("M3").Value2 = $snap[8].M
$ws.Range("O3").Value2 = $snap[8].O
$ws.Range("P3").Value2 = $snap[8].P

$ws.Range("D4").Value2 = $snap[17].D
$ws.Range("H4").Value2 = $snap[17].H
$ws.Range("I4").Value2 = $snap[17].I
$ws.Range("J4").Value2 = $snap[17].J
$ws.Range("K4").Value2 = $snap[17].K
$ws.Range("L4").Value2 = $snap[17].L
$ws.Range("M4").Value2 = $snap[17].M
$ws.Range("O4").Value2 = $snap[17].O
$ws.Range("P4").Value2 = $snap[17].P

$ws.Range("D5").Value2 = $snap[10].D
$ws.Range("H5").Value2 = $snap[10].H
$ws.Range("I5").Value2 = $snap[10].I
$ws.Range("J5").Value2 = $snap[10].J
$ws.Range("K5").Value2 = $snap[10].K
$ws.Range("L5").Value2 = $snap[10].L
$ws.Range("M5").Value2 = $snap[10].M
$ws.Range("O5").Value2 = $snap[10].O
$ws.Range("P5").Value2 = $snap[10].P

$ws.Range("D6").Value2 = $snap[16].D
$ws.Range("H6").Value2 = $snap[16].H
$ws.Range("I6").Value2 = $snap[16].I
$ws.Range("J6").Value2 = $snap[16].J
$ws.Range("K6").Value2 = $snap[16].K
$ws.Range("L6").Value2 = $snap[16].L
$ws.Range("M6").Value2 = $snap[16].M
$ws.Range("O6").Value2 = $snap[16].O
$ws.Range("P6").Value2 = $snap[16].P

$ws.Range("D8").Value2 = $snap[3].D
$ws.Range("H8").Value2 = $snap[3].H
$ws.Range("I8").Value2 = $snap[3].I
$ws.Range("J8").Value2 = $snap[3].J
$ws.Range("K8").Value2 = $snap[3].K
$ws.Range("L8").Value2 = $snap[3].L
$ws.Range("M8").Value2 = $snap[3].M
$ws.Range("O8").Value2 = $snap[3].O
$ws.Range("P8").Value2 = $snap[3].P

$ws.Range("D10").Value2 = $snap[15].D
$ws.Range("H10").Value2 = $snap[15].H
$ws.Range("I10").Value2 = $snap[15].I
$ws.Range("J10").Value2 = $snap[15].J
$ws.Range("K10").Value2 = $snap[15].K
$ws.Range("L10").Value2 = $snap[15].L
$ws.Range("M10").Value2 = $snap[15].M
$ws.Range("O10").Value2 = $snap[15].O
$ws.Range("P10").Value2 = $snap[15].P

$ws.Range("D11").Value2 = $snap[6].D
$ws.Range("H11").Value2 = $snap[6].H
$ws.Range("I11").Value2 = $snap[6].I
$ws.Range("J11").Value2 = $snap[6].J
$ws.Range("K11").Value2 = $snap[6].K
$ws.Range("L11").Value2 = $snap[6].L
$ws.Range("M11").Value2 = $snap[6].M
$ws.Range("O11").Value2 = $snap[6].O
$ws.Range("P11").Value2 = $snap[6].P

$ws.Range("D12").Value2 = $snap[2].D
$ws.Range("H12").Value2 = $snap[2].H
$ws.Range("I12").Value2 = $snap[2].I
$ws.Range("J12").Value2 = $snap[2].J
$ws.Range("K12").Value2 = $snap[2].K
$ws.Range("L12").Value2 = $snap[2].L
$ws.Range("M12").Value2 = $snap[2].M
$ws.Range("O12").Value2 = $snap[2].O
$ws.Range("P12").Value2 = $snap[2].P

$ws.Range("D13").Value2 = $snap[4].D
$ws.Range("H13").Value2 = $snap[4].H
$ws.Range("I13").Value2 = $snap[4].I
$ws.Range("J13").Value2 = $snap[4].J
$ws.Range("K13").Value2 = $snap[4].K
$ws.Range("L13").Value2 = $snap[4].L
$ws.Range("M13").Value2 = $snap[4].M
$ws.Range("O13").Value2 = $snap[4].O
$ws.Range("P13").Value2 = $snap[4].P

$ws.Range("D15").Value2 = $snap[11].D
$ws.Range("H15").Value2 = $snap[11].H
$ws.Range("I15").Value2 = $snap[11].I
$ws.Range("J15").Value2 = $snap[11].J
$ws.Range("K15").Value2 = $snap[11].K
$ws.Range("L15").Value2 = $snap[11].L
$ws.Range("M15").Value2 = $snap[11].M
$ws.Range("O15").Value2 = $snap[11].O
$ws.Range("P15").Value2 = $snap[11].P

$ws.Range("D16").Value2 = $snap[13].D
$ws.Range("H16").Value2 = $snap[13].H
$ws.Range("I16").Value2 = $snap[13].I
$ws.Range("J16").Value2 = $snap[13].J
$ws.Range("K16").Value2 = $snap[13].K
$ws.Range("L16").Value2 = $snap[13].L
$ws.Range("M16").Value2 = $snap[13].M
$ws.Range("O16").Value2 = $snap[13].O
$ws.Range("P16").Value2 = $snap[13].P

$ws.Range("D17").Value2 = $snap[12].D
$ws.Range("H17").Value2 = $snap[12].H
$ws.Range("I17").Value2 = $snap[12].I
$ws.Range("J17").Value2 = $snap[12].J
$ws.Range("K17").Value2 = $snap[12].K
$ws.Range("L17").Value2 = $snap[12].L
$ws.Range("M17").Value2 = $snap[12].M
$ws.Range("O17").Value2 = $snap[12].O
$ws.Range("P17").Value2 = $snap[12].P

